# Weekly price update: a new week of "Cebollín" (Mercado Mayorista Lo
# Valledor de Santiago) price data is inserted at the top of the existing
# data block (row 1057), pushing all the previously-recorded weeks down by
# six rows (one block = Extra/Primera/Segunda x2 origins).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (rows 1057:1099) down by 6 rows to make room for
# the new week's 6 records. Excel copies the row-2 style (date format on
# column D) down onto the freshly inserted rows automatically.
$ws.Rows("1057:1062").Insert()

# Columns that are constant for every record in this sheet.
$marketId   = 6
$market     = 'Mercado Mayorista Lo Valledor de Santiago'
$region     = 'Metropolitana'
$codreg     = 13
$catId      = 100112037
$categoria  = 'Cebollín'
$variedad   = 'Sin especificar'
$unidad     = '$/paquete 36 unidades'
$kgOuds     = 36
$clasif     = 'Hortaliza'

# New week's data: date 2021-11-09 (serial 44509).
$newRows = @(
    @{ Row=1057; Calidad='Extra';   Vol=910;  PMin=2200; PMax=2300; PProm=2248; Origen='Provincia de Chacabuco'; PKg=62 },
    @{ Row=1058; Calidad='Extra';   Vol=790;  PMin=2000; PMax=2200; PProm=2104; Origen='Región Metropolitana';   PKg=58 },
    @{ Row=1059; Calidad='Primera'; Vol=1020; PMin=1900; PMax=2000; PProm=1952; Origen='Provincia de Chacabuco'; PKg=54 },
    @{ Row=1060; Calidad='Primera'; Vol=990;  PMin=1800; PMax=1900; PProm=1856; Origen='Región Metropolitana';   PKg=52 },
    @{ Row=1061; Calidad='Segunda'; Vol=280;  PMin=1600; PMax=1600; PProm=1600; Origen='Provincia de Chacabuco'; PKg=44 },
    @{ Row=1062; Calidad='Segunda'; Vol=300;  PMin=1600; PMax=1600; PProm=1600; Origen='Región Metropolitana';   PKg=44 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $marketId
    $ws.Cells.Item($row, 2).Value  = $market
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = 44509
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Vol
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgOuds
    $ws.Cells.Item($row, 18).Value = $clasif
}
